$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("I2").Value = 165
$ws.Range("J2").Value = 3587.96

# Row 3
$ws.Range("B3").Value = "4201118395"
$ws.Range("C3").Value = "Molla Pharmacy"
$ws.Range("D3").Value = "Pangsha Bzaer Rajbari                                       "
$ws.Range("E3").Value = "KC82                                                        "
$ws.Range("F3").Value = "SKINV430-571195     "
$ws.Range("G3").Value = "03 Jan 2021"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 3379.15

# Row 4
$ws.Range("B4").Value = "4201118043"
$ws.Range("C4").Value = "The Medical Center"
$ws.Range("D4").Value = "In Fornt of Upzela Pangsha                                  "
$ws.Range("E4").Value = "KC81                                                        "
$ws.Range("F4").Value = "SKINV430-576011     "
$ws.Range("G4").Value = "13 Jan 2021"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 34
$ws.Range("J4").Value = 22984.02

# Row 5
$ws.Range("B5").Value = "90000006012"
$ws.Range("C5").Value = "Al Bab Medical Hall"
$ws.Range("D5").Value = "Kamalpur,Shampur,Daulotpur                                  "
$ws.Range("E5").Value = "KC32                                                        "
$ws.Range("F5").Value = "SKINV430-577547     "
$ws.Range("G5").Value = "17 Jan 2021"
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 29944.37

# Row 6
$ws.Range("B6").Value = "43061101800"
$ws.Range("C6").Value = "Renasase Med Centre"
$ws.Range("D6").Value = "Khalishakundi Mirpur                                        "
$ws.Range("E6").Value = "KC33                                                        "
$ws.Range("F6").Value = "SKINV430-577428     "
$ws.Range("G6").Value = "17 Jan 2021"
$ws.Range("I6").Value = 30
$ws.Range("J6").Value = 20040.05

# Row 7
$ws.Range("B7").Value = "90000008724"
$ws.Range("C7").Value = "Shefa Clinic"
$ws.Range("D7").Value = "Ershadpur, Alamdanga                                        "
$ws.Range("E7").Value = "KC63                                                        "
$ws.Range("F7").Value = "SKINV430-566768     "
$ws.Range("G7").Value = "21 Dec 2020"
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 27
$ws.Range("J7").Value = 12000.73

# Row 8
$ws.Range("B8").Value = "43061104555"
$ws.Range("C8").Value = "Opu Pharmacy"
$ws.Range("D8").Value = "Dorshona Bus Stand Dorshona                                 "
$ws.Range("E8").Value = "KC53                                                        "
$ws.Range("F8").Value = "SKINV430-580913     "
$ws.Range("G8").Value = "25 Jan 2021"
$ws.Range("I8").Value = 22
$ws.Range("J8").Value = 19458.12

# Row 9
$ws.Range("B9").Value = "43061109604"
$ws.Range("C9").Value = "Dr Eanur"
$ws.Range("D9").Value = "Sahapur Andulbaria Jibonnagar                               "
$ws.Range("E9").Value = "KC53                                                        "
$ws.Range("F9").Value = "SKINV430-580688     "
$ws.Range("G9").Value = "25 Jan 2021"
$ws.Range("I9").Value = 22
$ws.Range("J9").Value = 9839.27

# Row 10
$ws.Range("B10").Value = "90000047000"
$ws.Range("C10").Value = "Doctor Bari Pharmacy"
$ws.Range("D10").Value = "Badargonj Bazar.Choudanga                                   "
$ws.Range("E10").Value = "KC65                                                        "
$ws.Range("F10").Value = "SKINV430-581127     "
$ws.Range("G10").Value = "26 Jan 2021"
$ws.Range("I10").Value = 21
$ws.Range("J10").Value = 34716.14

# Row 11
$ws.Range("B11").Value = "43061104892"
$ws.Range("C11").Value = "Dr Belal Med Hall"
$ws.Range("D11").Value = "Allardorga Bazar ,Daulotpur                                 "
$ws.Range("E11").Value = "KC34                                                        "
$ws.Range("F11").Value = "SKINV430-581443     "
$ws.Range("G11").Value = "27 Jan 2021"
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 25198.49

# Row 12
$ws.Range("B12").Value = "43061110163"
$ws.Range("C12").Value = "Salam Medicine Corner"
$ws.Range("D12").Value = "Kurulgachi Bazar Dorshona Chuadanga                         "
$ws.Range("E12").Value = "KC53                                                        "
$ws.Range("F12").Value = "SKINV430-581951     "
$ws.Range("G12").Value = "28 Jan 2021"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 19
$ws.Range("J12").Value = 19741.02

# Row 13
$ws.Range("B13").Value = "90000045498"
$ws.Range("C13").Value = "Abdullah Pharmacy"
$ws.Range("D13").Value = "Halt Station Dorsona                                        "
$ws.Range("E13").Value = "KC53                                                        "
$ws.Range("F13").Value = "SKINV430-582184     "
$ws.Range("G13").Value = "28 Jan 2021"
$ws.Range("I13").Value = 19
$ws.Range("J13").Value = 9954.200000000001

# Row 14
$ws.Range("B14").Value = "90000003726"
$ws.Range("C14").Value = "M S Selina Drug House"
$ws.Range("D14").Value = "Fultola Chowrhas More                                       "
$ws.Range("E14").Value = "KC11                                                        "
$ws.Range("F14").Value = "SKINV430-582233     "
$ws.Range("G14").Value = "29 Jan 2021"
$ws.Range("I14").Value = 18
$ws.Range("J14").Value = 1500.95

# Row 15
$ws.Range("F15").Value = "SKINV430-582653     "
$ws.Range("G15").Value = "30 Jan 2021"
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 12761.05

# Row 16
$ws.Range("B16").Value = "90000029247"
$ws.Range("C16").Value = "Amena Medical Stor"
$ws.Range("D16").Value = "Allardarga Bazar.Doulatpur Kushtia.                         "
$ws.Range("E16").Value = "KC34                                                        "
$ws.Range("F16").Value = "SKINV430-583047     "
$ws.Range("G16").Value = "31 Jan 2021"
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 29871.3

# Row 17
$ws.Range("B17").Value = "90000044336"
$ws.Range("C17").Value = "Maa Fatema Medical Hall"
$ws.Range("D17").Value = "Nandolaur.Kumarkhali.Kushtia                                "
$ws.Range("E17").Value = "KC22                                                        "
$ws.Range("F17").Value = "SKINV430-583079     "
$ws.Range("G17").Value = "31 Jan 2021"
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 5139.4

# Row 18
$ws.Range("B18").Value = "43061101908"
$ws.Range("C18").Value = "Ashif Pharmacy"
$ws.Range("D18").Value = "Hospital Road Kushtia                                       "
$ws.Range("E18").Value = "KC24                                                        "
$ws.Range("F18").Value = "SKINV430-583086     "
$ws.Range("G18").Value = "31 Jan 2021"
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 4953.19

# Row 19
$ws.Range("B19").Value = "4201117970"
$ws.Range("C19").Value = "Bishwa Nath Medical "
$ws.Range("D19").Value = "Tatulia Bazar Baliakandi                                    "
$ws.Range("E19").Value = "KC83                                                        "
$ws.Range("F19").Value = "SKINV430-583463     "
$ws.Range("G19").Value = "01 Feb 2021"
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 19940.27

# Row 20
$ws.Range("B20").Value = "90000030489"
$ws.Range("C20").Value = "Emon Pharmacy"
$ws.Range("D20").Value = "Raghunatpur.Chuadanga                                       "
$ws.Range("E20").Value = "KC53                                                        "
$ws.Range("F20").Value = "SKINV430-583568     "
$ws.Range("G20").Value = "01 Feb 2021"
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 19541.16

# Row 21
$ws.Range("B21").Value = "43061000918"
$ws.Range("C21").Value = "Orient Pharmacy"
$ws.Range("D21").Value = "Circit House Meherpur                                       "
$ws.Range("E21").Value = "KC71                                                        "
$ws.Range("F21").Value = "SKINV430-570942     "
$ws.Range("G21").Value = "02 Jan 2021"
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = 15
$ws.Range("J21").Value = 14914.1
